# Natmi following Dr Hou advice
# Re-run of the NATMI edge-weight table for Lama2-Itgb1 (LR-pairs_lrc2p,
# YoungD0): "Sending cluster" / "Target cluster" now iterate over all
# three clusters (ECs, FAPs, sCs) instead of just FAPs/sCs, and the
# per-edge statistics (columns E:T) are recomputed accordingly. The
# table grows from 6 data rows (rows 2-7) to 9 data rows (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Lama2"
$ws.Cells.Item(2,3).Value = "Itgb1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 2.452389333333333
$ws.Cells.Item(2,8).Value = 7.357168
$ws.Cells.Item(2,9).Value = 0.007993767302975028
$ws.Cells.Item(2,10).Value = 0.007993767302975028
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 112.513392
$ws.Cells.Item(2,14).Value = 337.540176
$ws.Cells.Item(2,15).Value = 0.3275312977368564
$ws.Cells.Item(2,16).Value = 0.3275312977368564
$ws.Cells.Item(2,17).Value = 275.926642397952
$ws.Cells.Item(2,18).Value = 2483.339781581568
$ws.Cells.Item(2,19).Value = 0.002618208978549861
$ws.Cells.Item(2,20).Value = 0.002618208978549861

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Lama2"
$ws.Cells.Item(3,3).Value = "Itgb1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 2.452389333333333
$ws.Cells.Item(3,8).Value = 7.357168
$ws.Cells.Item(3,9).Value = 0.007993767302975028
$ws.Cells.Item(3,10).Value = 0.007993767302975028
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 106.314466
$ws.Cells.Item(3,14).Value = 318.943398
$ws.Cells.Item(3,15).Value = 0.3094859589441663
$ws.Cells.Item(3,16).Value = 0.3094859589441664
$ws.Cells.Item(3,17).Value = 260.7244623974293
$ws.Cells.Item(3,18).Value = 2346.520161576864
$ws.Cells.Item(3,19).Value = 0.002473958739337748
$ws.Cells.Item(3,20).Value = 0.002473958739337749

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Lama2"
$ws.Cells.Item(4,3).Value = "Itgb1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 2.452389333333333
$ws.Cells.Item(4,8).Value = 7.357168
$ws.Cells.Item(4,9).Value = 0.007993767302975028
$ws.Cells.Item(4,10).Value = 0.007993767302975028
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 124.6916553333333
$ws.Cells.Item(4,14).Value = 374.074966
$ws.Cells.Item(4,15).Value = 0.3629827433189773
$ws.Cells.Item(4,16).Value = 0.3629827433189773
$ws.Cells.Item(4,17).Value = 305.7924854951431
$ws.Cells.Item(4,18).Value = 2752.132369456288
$ws.Cells.Item(4,19).Value = 0.002901599585087418
$ws.Cells.Item(4,20).Value = 0.002901599585087418

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Lama2"
$ws.Cells.Item(5,3).Value = "Itgb1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 243.5672963333334
$ws.Cells.Item(5,8).Value = 730.7018890000001
$ws.Cells.Item(5,9).Value = 0.7939278902575405
$ws.Cells.Item(5,10).Value = 0.7939278902575405
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 112.513392
$ws.Cells.Item(5,14).Value = 337.540176
$ws.Cells.Item(5,15).Value = 0.3275312977368564
$ws.Cells.Item(5,16).Value = 0.3275312977368564
$ws.Cells.Item(5,17).Value = 27404.5826907325
$ws.Cells.Item(5,18).Value = 246641.2442165925
$ws.Cells.Item(5,19).Value = 0.2600362322055367
$ws.Cells.Item(5,20).Value = 0.2600362322055367

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Lama2"
$ws.Cells.Item(6,3).Value = "Itgb1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 243.5672963333334
$ws.Cells.Item(6,8).Value = 730.7018890000001
$ws.Cells.Item(6,9).Value = 0.7939278902575405
$ws.Cells.Item(6,10).Value = 0.7939278902575405
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 106.314466
$ws.Cells.Item(6,14).Value = 318.943398
$ws.Cells.Item(6,15).Value = 0.3094859589441663
$ws.Cells.Item(6,16).Value = 0.3094859589441664
$ws.Cells.Item(6,17).Value = 25894.72704474209
$ws.Cells.Item(6,18).Value = 233052.5434026788
$ws.Cells.Item(6,19).Value = 0.2457095344488738
$ws.Cells.Item(6,20).Value = 0.2457095344488738

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Lama2"
$ws.Cells.Item(7,3).Value = "Itgb1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 243.5672963333334
$ws.Cells.Item(7,8).Value = 730.7018890000001
$ws.Cells.Item(7,9).Value = 0.7939278902575405
$ws.Cells.Item(7,10).Value = 0.7939278902575405
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 124.6916553333333
$ws.Cells.Item(7,14).Value = 374.074966
$ws.Cells.Item(7,15).Value = 0.3629827433189773
$ws.Cells.Item(7,16).Value = 0.3629827433189773
$ws.Cells.Item(7,17).Value = 30370.80936486787
$ws.Cells.Item(7,18).Value = 273337.2842838108
$ws.Cells.Item(7,19).Value = 0.28818212360313
$ws.Cells.Item(7,20).Value = 0.2881821236031301

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Lama2"
$ws.Cells.Item(8,3).Value = "Itgb1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 60.76799533333334
$ws.Cells.Item(8,8).Value = 182.303986
$ws.Cells.Item(8,9).Value = 0.1980783424394845
$ws.Cells.Item(8,10).Value = 0.1980783424394845
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 112.513392
$ws.Cells.Item(8,14).Value = 337.540176
$ws.Cells.Item(8,15).Value = 0.3275312977368564
$ws.Cells.Item(8,16).Value = 0.3275312977368564
$ws.Cells.Item(8,17).Value = 6837.213279993504
$ws.Cells.Item(8,18).Value = 61534.91951994153
$ws.Cells.Item(8,19).Value = 0.06487685655276978
$ws.Cells.Item(8,20).Value = 0.06487685655276978

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Lama2"
$ws.Cells.Item(9,3).Value = "Itgb1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 60.76799533333334
$ws.Cells.Item(9,8).Value = 182.303986
$ws.Cells.Item(9,9).Value = 0.1980783424394845
$ws.Cells.Item(9,10).Value = 0.1980783424394845
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 106.314466
$ws.Cells.Item(9,14).Value = 318.943398
$ws.Cells.Item(9,15).Value = 0.3094859589441663
$ws.Cells.Item(9,16).Value = 0.3094859589441664
$ws.Cells.Item(9,17).Value = 6460.516973753825
$ws.Cells.Item(9,18).Value = 58144.65276378443
$ws.Cells.Item(9,19).Value = 0.06130246575595481
$ws.Cells.Item(9,20).Value = 0.06130246575595483

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Lama2"
$ws.Cells.Item(10,3).Value = "Itgb1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 60.76799533333334
$ws.Cells.Item(10,8).Value = 182.303986
$ws.Cells.Item(10,9).Value = 0.1980783424394845
$ws.Cells.Item(10,10).Value = 0.1980783424394845
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 124.6916553333333
$ws.Cells.Item(10,14).Value = 374.074966
$ws.Cells.Item(10,15).Value = 0.3629827433189773
$ws.Cells.Item(10,16).Value = 0.3629827433189773
$ws.Cells.Item(10,17).Value = 7577.261929401609
$ws.Cells.Item(10,18).Value = 68195.35736461448
$ws.Cells.Item(10,19).Value = 0.07189902013075988
$ws.Cells.Item(10,20).Value = 0.0718990201307599
